$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.064.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.284.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +4.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.283.55"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0998"
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.691.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.009.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.264.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "313.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.169"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E28").Value = "  -3.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("E30").Value = "  -1.75%  "
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.07"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.379"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  -1.27%  "
$ws.Range("E39").Value = "  -1.08%  "
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "139.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "285.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.62%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("E44").Value = "  +1.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0494"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("E46").Value = "  +0.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.02%  "
$ws.Range("E48").Value = "  -1.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.40%  "
